$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

$ws.Range("B38").Value = "SingleUseId40"
$ws.Range("C38").Value = "Default"
$ws.Range("D38").Value = "Center"
$ws.Range("E38").Value = "LTR"
$ws.Range("F38").Value = "<value> ms"
$ws.Range("B39").Value = "SingleUseId42"
$ws.Range("C39").Value = "Default"
$ws.Range("D39").Value = "Center"
$ws.Range("E39").Value = "LTR"
$ws.Range("F39").Value = "<value>"
$ws.Range("B40").Value = "SingleUseId43"
$ws.Range("C40").Value = "Default"
$ws.Range("D40").Value = "Left"
$ws.Range("E40").Value = "LTR"
$ws.Range("F40").Value = "'0"
$ws.Range("B41").Value = "SingleUseId44"
$ws.Range("C41").Value = "Default"
$ws.Range("D41").Value = "Left"
$ws.Range("E41").Value = "LTR"
$ws.Range("F41").Value = "'0"
$ws.Range("B42").Value = "SingleUseId45"
$ws.Range("C42").Value = "Default"
$ws.Range("D42").Value = "Left"
$ws.Range("E42").Value = "LTR"
$ws.Range("F42").Value = "Stamps Number"
$ws.Range("B43").Value = "SingleUseId46"
$ws.Range("C43").Value = "Large"
$ws.Range("D43").Value = "Left"
$ws.Range("E43").Value = "LTR"
$ws.Range("F43").Value = "X"
$ws.Range("B44").Value = "SingleUseId47"
$ws.Range("C44").Value = "Default"
$ws.Range("D44").Value = "Left"
$ws.Range("E44").Value = "LTR"
$ws.Range("F44").Value = "Repeat"
$ws.Range("B45").Value = "SingleUseId48"
$ws.Range("C45").Value = "Default"
$ws.Range("D45").Value = "Center"
$ws.Range("E45").Value = "LTR"
$ws.Range("F45").Value = "<value>"
$ws.Range("B46").Value = "SingleUseId49"
$ws.Range("C46").Value = "Default"
$ws.Range("D46").Value = "Left"
$ws.Range("E46").Value = "LTR"
$ws.Range("F46").Value = "'0"
$ws.Range("B47").Value = "SingleUseId52"
$ws.Range("C47").Value = "Default"
$ws.Range("D47").Value = "Left"
$ws.Range("E47").Value = "LTR"
$ws.Range("F47").Value = "Single"
$ws.Range("B48").Value = "SingleUseId53"
$ws.Range("C48").Value = "Default"
$ws.Range("D48").Value = "Left"
$ws.Range("E48").Value = "LTR"
$ws.Range("F48").Value = "Continuous"
$ws.Range("B49").Value = "SingleUseId54"
$ws.Range("C49").Value = "Default"
$ws.Range("D49").Value = "Left"
$ws.Range("E49").Value = "LTR"
$ws.Range("F49").Value = "Value from main.c: <value>"
$ws.Range("B50").Value = "SingleUseId55"
$ws.Range("C50").Value = "Default"
$ws.Range("D50").Value = "Left"
$ws.Range("E50").Value = "LTR"
$ws.Range("F50").Value = "THRESHOLD"
$ws.Range("B51").Value = "SingleUseId56"
$ws.Range("C51").Value = "Default"
$ws.Range("D51").Value = "Left"
$ws.Range("E51").Value = "LTR"
$ws.Range("F51").Value = "SLOPE"
$ws.Range("B52").Value = "SingleUseId57"
$ws.Range("C52").Value = "Default"
$ws.Range("D52").Value = "Center"
$ws.Range("E52").Value = "LTR"
$ws.Range("F52").Value = "Detect"
$ws.Range("B53").Value = "SingleUseId58"
$ws.Range("C53").Value = "Default"
$ws.Range("D53").Value = "Center"
$ws.Range("E53").Value = "LTR"
$ws.Range("F53").Value = "<value> mV"
$ws.Range("B54").Value = "SingleUseId59"
$ws.Range("C54").Value = "Default"
$ws.Range("D54").Value = "Left"
$ws.Range("E54").Value = "LTR"
$ws.Range("F54").Value = "'0"
$ws.Range("B55").Value = "SingleUseId62"
$ws.Range("C55").Value = "Default"
$ws.Range("D55").Value = "Left"
$ws.Range("E55").Value = "LTR"
$ws.Range("F55").Value = "External"
$ws.Range("B56").Value = "SingleUseId63"
$ws.Range("C56").Value = "Default"
$ws.Range("D56").Value = "Left"
$ws.Range("E56").Value = "LTR"
$ws.Range("F56").Value = "Internal`nRubid"
$ws.Range("B57").Value = "SingleUseId64"
$ws.Range("C57").Value = "Default"
$ws.Range("D57").Value = "Center"
$ws.Range("E57").Value = "LTR"
$ws.Range("F57").Value = "Internal`nQuartz"
$ws.Range("B58").Value = "SingleUseId65"
$ws.Range("C58").Value = "Default"
$ws.Range("D58").Value = "Center"
$ws.Range("E58").Value = "LTR"
$ws.Range("F58").Value = "<value>"
$ws.Range("B59").Value = "SingleUseId66"
$ws.Range("C59").Value = "Default"
$ws.Range("D59").Value = "Left"
$ws.Range("E59").Value = "LTR"
$ws.Range("F59").Value = "HF INPUT"
$ws.Range("B60").Value = "SingleUseId67"
$ws.Range("C60").Value = "Large"
$ws.Range("D60").Value = "Left"
$ws.Range("E60").Value = "LTR"
$ws.Range("F60").Value = "ON"
$ws.Range("B61").Value = "SingleUseId68"
$ws.Range("C61").Value = "Large"
$ws.Range("D61").Value = "Left"
$ws.Range("E61").Value = "LTR"
$ws.Range("F61").Value = "OFF"
$ws.Range("B62").Value = "SingleUseId69"
$ws.Range("C62").Value = "Default"
$ws.Range("D62").Value = "Left"
$ws.Range("E62").Value = "LTR"
$ws.Range("F62").Value = "GATE"
$ws.Range("B63").Value = "SingleUseId70"
$ws.Range("C63").Value = "Default"
$ws.Range("D63").Value = "Center"
$ws.Range("E63").Value = "LTR"
$ws.Range("F63").Value = "<value> ms"
$ws.Range("B64").Value = "SingleUseId71"
$ws.Range("C64").Value = "Default"
$ws.Range("D64").Value = "Left"
$ws.Range("E64").Value = "LTR"
$ws.Range("F64").Value = "MES SETUP"
$ws.Range("B65").Value = "SingleUseId73"
$ws.Range("C65").Value = "Default"
$ws.Range("D65").Value = "Center"
$ws.Range("E65").Value = "LTR"
$ws.Range("F65").Value = "<value>"
$ws.Range("B66").Value = "SingleUseId74"
$ws.Range("C66").Value = "Default"
$ws.Range("D66").Value = "Left"
$ws.Range("E66").Value = "LTR"
$ws.Range("F66").Value = "'0"
$ws.Range("B67").Value = "SingleUseId75"
$ws.Range("C67").Value = "Default"
$ws.Range("D67").Value = "Left"
$ws.Range("E67").Value = "LTR"
$ws.Range("F67").Value = "Stamps Number"
$ws.Range("B68").Value = "SingleUseId76"
$ws.Range("C68").Value = "Large"
$ws.Range("D68").Value = "Left"
$ws.Range("E68").Value = "LTR"
$ws.Range("F68").Value = "X"
$ws.Range("B69").Value = "SingleUseId77"
$ws.Range("C69").Value = "Default"
$ws.Range("D69").Value = "Left"
$ws.Range("E69").Value = "LTR"
$ws.Range("F69").Value = "Repeat"
$ws.Range("B70").Value = "SingleUseId78"
$ws.Range("C70").Value = "Default"
$ws.Range("D70").Value = "Center"
$ws.Range("E70").Value = "LTR"
$ws.Range("F70").Value = "<value>"
$ws.Range("B71").Value = "SingleUseId79"
$ws.Range("C71").Value = "Default"
$ws.Range("D71").Value = "Left"
$ws.Range("E71").Value = "LTR"
$ws.Range("F71").Value = "'0"
$ws.Range("B72").Value = "SingleUseId82"
$ws.Range("C72").Value = "Default"
$ws.Range("D72").Value = "Left"
$ws.Range("E72").Value = "LTR"
$ws.Range("F72").Value = "Single"
$ws.Range("B73").Value = "SingleUseId83"
$ws.Range("C73").Value = "Default"
$ws.Range("D73").Value = "Left"
$ws.Range("E73").Value = "LTR"
$ws.Range("F73").Value = "Continuous"
$ws.Range("B74").Value = "SingleUseId84"
$ws.Range("C74").Value = "Default"
$ws.Range("D74").Value = "Left"
$ws.Range("E74").Value = "LTR"
$ws.Range("F74").Value = "Value from main.c: <value>"
$ws.Range("B75").Value = "SingleUseId85"
$ws.Range("C75").Value = "Default"
$ws.Range("D75").Value = "Left"
$ws.Range("E75").Value = "LTR"
$ws.Range("F75").Value = "Value from main.c: <value>"
# Rows 76-81 previously held SingleUseId80..85 (now shifted up into 38-75); clear their old content
$ws.Range("B76:F81").ClearContents()
